$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H9").Value = 6031.4707
$ws_ALC.Range("I9").Value = 8504.333000000001
$ws_ALC.Range("J9").Value = 96.59999999999999
$ws_ALC.Range("K9").Value = 8504.333000000001
$ws_ALC.Range("L9").Value = 96.59999999999999
$ws_ALC.Range("M9").Value = -8335.333000000001
$ws_ALC.Range("N9").Value = -434.6
$ws_ALC.Range("H40").Value = 3589089.5
$ws_ALC.Range("J40").Value = 6178980
$ws_ALC.Range("L40").Value = 6178980
$ws_ALC.Range("N40").Value = -6179330
$ws_ALC.Range("H42").Value = 1562.6
$ws_ALC.Range("I42").Value = 86.59999999999999
$ws_ALC.Range("J42").Value = 4514.6
$ws_ALC.Range("K42").Value = 259.8
$ws_ALC.Range("L42").Value = 13543.8
$ws_ALC.Range("M42").Value = -29.79999999999995
$ws_ALC.Range("N42").Value = -14003.8
$ws_ALC.Range("H62").Value = 1259.5555
$ws_ALC.Range("I62").Value = 1233.8572
$ws_ALC.Range("K62").Value = 1233.8572
$ws_ALC.Range("M62").Value = -609.8571999999999
$ws_ALC.Range("H65").Value = 1259.5555
$ws_ALC.Range("I65").Value = 1233.8572
$ws_ALC.Range("K65").Value = 6169.286
$ws_ALC.Range("M65").Value = -3049.286
$ws_ALC.Range("H112").Value = 78801.80499999999
$ws_ALC.Range("I112").Value = 112172.11
$ws_ALC.Range("J112").Value = 61135.176
$ws_ALC.Range("K112").Value = 336516.33
$ws_ALC.Range("L112").Value = 183405.528
$ws_ALC.Range("M112").Value = -335408.33
$ws_ALC.Range("N112").Value = -185621.528
$ws_ALC.Range("H116").Value = 7625.6665
$ws_ALC.Range("I116").Value = 7625.6665
$ws_ALC.Range("K116").Value = 7625.6665
$ws_ALC.Range("M116").Value = -4183.6665
$ws_ALC.Range("H130").Value = 100000
$ws_ALC.Range("J130").Value = 100000
$ws_ALC.Range("L130").Value = 100000
$ws_ALC.Range("N130").Value = -110040
$ws_ALC.Range("H137").Value = 2778.4375
$ws_ALC.Range("I137").Value = 1896.7778
$ws_ALC.Range("K137").Value = 5690.3334
$ws_ALC.Range("M137").Value = -3140.3334
$ws_ALC.Range("H138").Value = 3539.6667
$ws_ALC.Range("J138").Value = 4412.9
$ws_ALC.Range("L138").Value = 13238.7
$ws_ALC.Range("N138").Value = -23518.7
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H4").Value = 250
$ws_ARM.Range("I4").Value = 200
$ws_ARM.Range("K4").Value = 200
$ws_ARM.Range("M4").Value = -84
$ws_ARM.Range("H61").Value = 52633948
$ws_ARM.Range("I61").Value = 90910456
$ws_ARM.Range("J61").Value = 3749.375
$ws_ARM.Range("K61").Value = 90910456
$ws_ARM.Range("L61").Value = 3749.375
$ws_ARM.Range("M61").Value = -90910244
$ws_ARM.Range("N61").Value = -4173.375
$ws_ARM.Range("H74").Value = 45457156
$ws_ARM.Range("I74").Value = 50001868
$ws_ARM.Range("K74").Value = 50001868
$ws_ARM.Range("M74").Value = -50000994
$ws_ARM.Range("H77").Value = 45457156
$ws_ARM.Range("I77").Value = 50001868
$ws_ARM.Range("K77").Value = 250009340
$ws_ARM.Range("M77").Value = -250004972
$ws_ARM.Range("H132").Value = 8337243.5
$ws_ARM.Range("I132").Value = 12501928
$ws_ARM.Range("K132").Value = 37505784
$ws_ARM.Range("M132").Value = -37503254
$ws_ARM.Range("H136").Value = 52633948
$ws_ARM.Range("I136").Value = 90910456
$ws_ARM.Range("J136").Value = 3749.375
$ws_ARM.Range("K136").Value = 272731368
$ws_ARM.Range("L136").Value = 11248.125
$ws_ARM.Range("M136").Value = -272728818
$ws_ARM.Range("N136").Value = -16348.125
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H107").Value = 168230.67
$ws_BSM.Range("I107").Value = 1888.3334
$ws_BSM.Range("K107").Value = 1888.3334
$ws_BSM.Range("M107").Value = 31.66660000000002
$ws_BSM.Range("H108").Value = 69999
$ws_BSM.Range("J108").Value = 69999
$ws_BSM.Range("L108").Value = 69999
$ws_BSM.Range("N108").Value = -77679
$ws_BSM.Range("H109").Value = 66666
$ws_BSM.Range("J109").Value = 66666
$ws_BSM.Range("L109").Value = 66666
$ws_BSM.Range("N109").Value = -69440
$ws_BSM.Range("H111").Value = 0
$ws_BSM.Range("J111").Value = 0
$ws_BSM.Range("L111").Value = 0
$ws_BSM.Range("N111").ClearContents()
$ws_BSM.Range("H112").Value = 66999
$ws_BSM.Range("J112").Value = 66999
$ws_BSM.Range("L112").Value = 66999
$ws_BSM.Range("N112").Value = -69953
$ws_BSM.Range("H118").Value = 187999.5
$ws_BSM.Range("J118").Value = 187999.5
$ws_BSM.Range("L118").Value = 187999.5
$ws_BSM.Range("N118").Value = -191313.5
$ws_BSM.Range("H122").Value = 49900
$ws_BSM.Range("J122").Value = 49900
$ws_BSM.Range("L122").Value = 49900
$ws_BSM.Range("N122").Value = -59700
$ws_BSM.Range("H134").Value = 21740976
$ws_BSM.Range("J134").Value = 3250
$ws_BSM.Range("L134").Value = 9750
$ws_BSM.Range("N134").Value = -14820
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 4572.1333
$ws_CRP.Range("I31").Value = 4794.727
$ws_CRP.Range("K31").Value = 4794.727
$ws_CRP.Range("M31").Value = -4499.727
$ws_CRP.Range("H34").Value = 4572.1333
$ws_CRP.Range("I34").Value = 4794.727
$ws_CRP.Range("K34").Value = 4794.727
$ws_CRP.Range("M34").Value = -4592.727
$ws_CRP.Range("H58").Value = 12504395
$ws_CRP.Range("I58").Value = 25006186
$ws_CRP.Range("J58").Value = 2603.55
$ws_CRP.Range("K58").Value = 25006186
$ws_CRP.Range("L58").Value = 2603.55
$ws_CRP.Range("M58").Value = -25005983
$ws_CRP.Range("N58").Value = -3009.55
$ws_CRP.Range("H112").Value = 29999
$ws_CRP.Range("J112").Value = 29999
$ws_CRP.Range("L112").Value = 29999
$ws_CRP.Range("N112").Value = -32953
$ws_CRP.Range("H136").Value = 12504395
$ws_CRP.Range("I136").Value = 25006186
$ws_CRP.Range("J136").Value = 2603.55
$ws_CRP.Range("K136").Value = 75018558
$ws_CRP.Range("L136").Value = 7810.650000000001
$ws_CRP.Range("M136").Value = -75016008
$ws_CRP.Range("N136").Value = -12910.65
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 45.857143
$ws_CUL.Range("J2").Value = 67.666664
$ws_CUL.Range("L2").Value = 405.999984
$ws_CUL.Range("N2").Value = -631.999984
$ws_CUL.Range("H109").Value = 1618.7
$ws_CUL.Range("I109").Value = 1465.2222
$ws_CUL.Range("J109").Value = 3000
$ws_CUL.Range("K109").Value = 4395.6666
$ws_CUL.Range("L109").Value = 9000
$ws_CUL.Range("M109").Value = -3355.6666
$ws_CUL.Range("N109").Value = -11080
$ws_CUL.Range("H122").Value = 1380
$ws_CUL.Range("I122").Value = 1250
$ws_CUL.Range("K122").Value = 11250
$ws_CUL.Range("M122").Value = -8800
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 2313.7144
$ws_GSM.Range("I80").Value = 2232.3333
$ws_GSM.Range("J80").Value = 2374.75
$ws_GSM.Range("K80").Value = 2232.3333
$ws_GSM.Range("L80").Value = 2374.75
$ws_GSM.Range("M80").Value = -1234.3333
$ws_GSM.Range("N80").Value = -4370.75
$ws_GSM.Range("H83").Value = 2313.7144
$ws_GSM.Range("I83").Value = 2232.3333
$ws_GSM.Range("J83").Value = 2374.75
$ws_GSM.Range("K83").Value = 11161.6665
$ws_GSM.Range("L83").Value = 11873.75
$ws_GSM.Range("M83").Value = -6169.666499999999
$ws_GSM.Range("N83").Value = -21857.75
$ws_GSM.Range("H126").Value = 3666.5
$ws_GSM.Range("I126").Value = 2566.6667
$ws_GSM.Range("J126").Value = 6966
$ws_GSM.Range("K126").Value = 7700.000100000001
$ws_GSM.Range("L126").Value = 20898
$ws_GSM.Range("M126").Value = -5230.000100000001
$ws_GSM.Range("N126").Value = -25838
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H100").Value = 9507029
$ws_LTW.Range("I100").Value = 10507374
$ws_LTW.Range("K100").Value = 10507374
$ws_LTW.Range("M100").Value = -10506833
$ws_LTW.Range("H132").Value = 20003406
$ws_LTW.Range("I132").Value = 20872814
$ws_LTW.Range("K132").Value = 62618442
$ws_LTW.Range("M132").Value = -62615912
$ws_LTW.Range("H136").Value = 1441.625
$ws_LTW.Range("I136").Value = 1405.0667
$ws_LTW.Range("K136").Value = 4215.2001
$ws_LTW.Range("M136").Value = -1665.2001
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H136").Value = 15626559
$ws_WVR.Range("I136").Value = 17242912
$ws_WVR.Range("K136").Value = 51728736
$ws_WVR.Range("M136").Value = -51726186
